$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force number-like Price cells to remain text (matches source inlineStr formatting)
# so values like '1.00' / '0.605' aren't silently coerced into floats.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D10", "D11", "D12", "D13", "D18", "D19", "D20", "D23", "D27", "D28", "D29", "D33", "D34", "D36", "D37", "D41", "D42", "D43", "D45", "D46", "D47", "D49", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) columns with latest scraped values
$ws.Range("D2").Value = '63.317.51'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '2.565.18'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '584.51'
$ws.Range("E5").Value = '  +3.24%  '
$ws.Range("D6").Value = '148.00'
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.605'
$ws.Range("E8").Value = '  +4.60%  '
$ws.Range("E9").Value = '  +4.42%  '
$ws.Range("D10").Value = '5.67'
$ws.Range("E10").Value = '  +1.41%  '
$ws.Range("D11").Value = '0.153'
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").Value = '0.358'
$ws.Range("E12").Value = '  +2.05%  '
$ws.Range("D13").Value = '27.49'
$ws.Range("E13").Value = '  +2.28%  '
$ws.Range("D14").Value = '3.027.36'
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").Value = '63.246.06'
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("E16").Value = '  +5.40%  '
$ws.Range("D17").Value = '2.584.24'
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("D18").Value = '11.37'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").Value = '342.91'
$ws.Range("E19").Value = '  +2.90%  '
$ws.Range("D20").Value = '4.43'
$ws.Range("E20").Value = '  +3.86%  '
$ws.Range("E21").Value = '  +1.78%  '
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").Value = '66.81'
$ws.Range("E23").Value = '  +3.44%  '
$ws.Range("D24").Value = '2.685.41'
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("E25").Value = '  +3.84%  '
$ws.Range("E26").Value = '  +0.86%  '
$ws.Range("D27").Value = '8.21'
$ws.Range("E27").Value = '  +13.59%  '
$ws.Range("D28").Value = '8.55'
$ws.Range("E28").Value = '  +3.15%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("E31").Value = '  +7.89%  '
$ws.Range("D32").Value = '0.0₃0828'
$ws.Range("E32").Value = '  +2.67%  '
$ws.Range("D33").Value = '466.80'
$ws.Range("E33").Value = '  +15.54%  '
$ws.Range("D34").Value = '176.92'
$ws.Range("E34").Value = '  +0.16%  '
$ws.Range("E35").Value = '  +3.52%  '
$ws.Range("D36").Value = '0.408'
$ws.Range("E36").Value = '  +3.36%  '
$ws.Range("D37").Value = '19.29'
$ws.Range("E37").Value = '  +2.71%  '
$ws.Range("E38").Value = '  +5.16%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = '151.67'
$ws.Range("E42").Value = '  +0.23%  '
$ws.Range("D43").Value = '3.82'
$ws.Range("E43").Value = '  +2.70%  '
$ws.Range("E44").Value = '  +3.44%  '
$ws.Range("D45").Value = '0.0554'
$ws.Range("E45").Value = '  +7.30%  '
$ws.Range("D46").Value = '0.615'
$ws.Range("E46").Value = '  +2.74%  '
$ws.Range("D47").Value = '0.0984'
$ws.Range("E47").Value = '  +3.03%  '
$ws.Range("E48").Value = '  +2.62%  '
$ws.Range("D49").Value = '18.52'
$ws.Range("E49").Value = '  +1.61%  '
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").Value = '11.39'
$ws.Range("E51").Value = '  -0.20%  '

# Drop the transient Text-format stamp so unaffected cell styling is untouched
foreach ($c in $textCells) {
    $ws.Range($c).ClearFormats()
}
